$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$targetCell = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $label = $t.Cell($i, 1).Range.Text
    if ($label -like "Micro results*") {
        $targetCell = $t.Cell($i, 2)
        break
    }
}

if ($targetCell -eq $null) {
    throw "Could not find 'Micro results' row"
}

$r = $targetCell.Range
$r.Text = ""

$r2 = $targetCell.Range
$r2.Collapse(0)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>18/08 – BLC – PERIPHERAL–LEFT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>18/08 – BLC – RED PORT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>16/08 – CRYPTO GIARDIA DIRECT ANTIGEN – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Giardia and Cryptosporidium not detected.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>16/08 – Salmonella PCR – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Salmonella not detected by PCR.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>16/08 – Campylobacter PCR – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Campylobacter not detected.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>16/08 – E.coli O157 &amp; Verotoxin PCR – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: E.coli O157/VTEC not detected.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>16/08 – C difficile(GDH) – Negative</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>16/08 – Shigella PCR – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Shigella species not detected.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>15/08 – BLC – WHITE PORT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>15/08 – BLC – PERIPHERAL–LEFT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>15/08 – BLC – RED PORT NO GROWTH AFTER 5 DAYS</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>15/08 – UC – MID STREAM URINE (MSU) NO SIGNIFICANT GROWTH</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>15/08 – UC – **No clear Result**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: Antibiotic usage advice, no microbiological findings.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>15/08 – RESPIRATORY PCR – **Negative**</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">  Summary: No microorganisms detected by BIOFIRE PCR.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>--------Previous result (1 year)--------</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>24/07 – EBV VCA IgG – Positive</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">19/06 – RESPIRATORY PCR – **Positive**  </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>**Summary:** Human Rhinovirus/Enterovirus detected.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>12/06 – EBV VCA IgG – Positive</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">20/02 – RESPIRATORY PCR – **Positive**  </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>**Summary:** Respiratory Syncytial Virus detected.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t xml:space="preserve">04/02 – RESPIRATORY PCR – **Positive**  </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:color w:val="0000FF"/><w:sz w:val="20"/></w:rPr><w:t>**Summary:** Respiratory Syncytial Virus detected.</w:t></w:r></w:p>'
$r2.InsertXML($xml)

Write-Host "Micro results cell updated."
